$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns (row 1) to English machine-readable names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case the Spanish connector words ("de", "del", "la", "las",
#        "los", "el", "y") inside the state/municipality name columns.
#        Every occurrence in the data is surrounded by plain spaces, so a
#        simple space-padded substring replace is equivalent to a whole
#        word replace here and is order independent. ---
$nameRange = $ws.Range("A1:B2341")

# xlPart = 2 (LookAt), xlByRows = 1 (SearchOrder) ; MatchCase = $true
$nameRange.Replace(" de ", " De ", 2, 1, $false, $false, $true, $true) | Out-Null
$nameRange.Replace(" del ", " Del ", 2, 1, $false, $false, $true, $true) | Out-Null
$nameRange.Replace(" la ", " La ", 2, 1, $false, $false, $true, $true) | Out-Null
$nameRange.Replace(" las ", " Las ", 2, 1, $false, $false, $true, $true) | Out-Null
$nameRange.Replace(" los ", " Los ", 2, 1, $false, $false, $true, $true) | Out-Null
$nameRange.Replace(" el ", " El ", 2, 1, $false, $false, $true, $true) | Out-Null
$nameRange.Replace(" y ", " Y ", 2, 1, $false, $false, $true, $true) | Out-Null

# --- 3. Remove the trailing footnote rows (the blank row 2336 plus the
#        source/notes rows 2337-2341), shrinking the used range down to
#        A1:D2335. ---
$ws.Range("A2336:A2341").EntireRow.Delete() | Out-Null
